$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '29.480.11'
Set-TextValue $ws.Range("E2") '  -0.14%  '
Set-TextValue $ws.Range("D3") '1.917.72'
Set-TextValue $ws.Range("E3") '  +0.52%  '
Set-TextValue $ws.Range("D5") '325.37'
Set-TextValue $ws.Range("E5") '  -0.17%  '
Set-TextValue $ws.Range("D6") '1.012'
Set-TextValue $ws.Range("E6") '  +0.62%  '
Set-TextValue $ws.Range("D7") '0.4802'
Set-TextValue $ws.Range("E7") '  -0.88%  '
Set-TextValue $ws.Range("D8") '0.4050'
Set-TextValue $ws.Range("E8") '  -0.59%  '
Set-TextValue $ws.Range("D9") '0.08207'
Set-TextValue $ws.Range("E9") '  +0.66%  '
Set-TextValue $ws.Range("D10") '1.008'
Set-TextValue $ws.Range("E10") '  -0.23%  '
Set-TextValue $ws.Range("D11") '23.40'
Set-TextValue $ws.Range("E11") '  -0.20%  '
Set-TextValue $ws.Range("D12") '1.922.58'
Set-TextValue $ws.Range("E12") '  +3.80%  '
Set-TextValue $ws.Range("D13") '6.047'
Set-TextValue $ws.Range("E13") '  +0.28%  '
Set-TextValue $ws.Range("D14") '7.216'
Set-TextValue $ws.Range("E14") '  +1.76%  '
Set-TextValue $ws.Range("D15") '91.28'
Set-TextValue $ws.Range("E15") '  +0.94%  '
Set-TextValue $ws.Range("D16") '0.06850'
Set-TextValue $ws.Range("E16") '  +1.45%  '
Set-TextValue $ws.Range("E17") '  +0.60%  '
Set-TextValue $ws.Range("E18") '  -0.65%  '
Set-TextValue $ws.Range("D19") '17.51'
Set-TextValue $ws.Range("E19") '  -1.10%  '
Set-TextValue $ws.Range("E20") '  +0.61%  '
Set-TextValue $ws.Range("D21") '29.482.68'
Set-TextValue $ws.Range("E21") '  -0.16%  '
Set-TextValue $ws.Range("D22") '5.661'
Set-TextValue $ws.Range("E22") '  +1.29%  '
Set-TextValue $ws.Range("E23") '  +0.39%  '
Set-TextValue $ws.Range("E24") '  +1.28%  '
Set-TextValue $ws.Range("D25") '2.147.88'
Set-TextValue $ws.Range("E25") '  +2.97%  '
Set-TextValue $ws.Range("B26") 'Monero'
Set-TextValue $ws.Range("C26") 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D26") '156.10'
Set-TextValue $ws.Range("E26") '  +1.22%  '
Set-TextValue $ws.Range("B27") 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range("C27") 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D27") '6.474'
Set-TextValue $ws.Range("E27") '  +2.66%  '
Set-TextValue $ws.Range("D28") '20.01'
Set-TextValue $ws.Range("E28") '  -0.31%  '
Set-TextValue $ws.Range("E29") '  -0.51%  '
Set-TextValue $ws.Range("D30") '120.44'
Set-TextValue $ws.Range("E30") '  +1.16%  '
Set-TextValue $ws.Range("D31") '1.011'
Set-TextValue $ws.Range("E31") '  -2.70%  '
Set-TextValue $ws.Range("D32") '0.09608'
Set-TextValue $ws.Range("E32") '  +0.51%  '
Set-TextValue $ws.Range("D33") '5.613'
Set-TextValue $ws.Range("E33") '  +1.36%  '
Set-TextValue $ws.Range("D34") '3.559'
Set-TextValue $ws.Range("E34") '  +0.12%  '
Set-TextValue $ws.Range("D35") '1.368'
Set-TextValue $ws.Range("E35") '  -1.92%  '
Set-TextValue $ws.Range("D36") '0.06301'
Set-TextValue $ws.Range("E36") '  +3.01%  '
Set-TextValue $ws.Range("D37") '0.02277'
Set-TextValue $ws.Range("E37") '  +0.45%  '
Set-TextValue $ws.Range("D38") '1.181'
Set-TextValue $ws.Range("E38") '  +0.56%  '
Set-TextValue $ws.Range("D39") '0.5925'
Set-TextValue $ws.Range("E39") '  -0.74%  '
Set-TextValue $ws.Range("E40") '  +2.27%  '
Set-TextValue $ws.Range("D41") '7.861'
Set-TextValue $ws.Range("E41") '  -0.87%  '
Set-TextValue $ws.Range("E42") '  -0.58%  '
Set-TextValue $ws.Range("E43") '  +0.07%  '
Set-TextValue $ws.Range("D44") '2.386'
Set-TextValue $ws.Range("E44") '  -2.11%  '
Set-TextValue $ws.Range("D45") '12.41'
Set-TextValue $ws.Range("E45") '  -0.52%  '
Set-TextValue $ws.Range("D46") '0.07470'
Set-TextValue $ws.Range("E46") '  -3.30%  '
Set-TextValue $ws.Range("D47") '0.5549'
Set-TextValue $ws.Range("E47") '  -0.41%  '
Set-TextValue $ws.Range("D48") '1.931'
Set-TextValue $ws.Range("E48") '  -1.51%  '
Set-TextValue $ws.Range("D49") '117.83'
Set-TextValue $ws.Range("E49") '  +2.49%  '
Set-TextValue $ws.Range("E50") '  +3.04%  '
Set-TextValue $ws.Range("D51") '71.88'
Set-TextValue $ws.Range("E51") '  -1.05%  '
